$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object "object[,]" 24,1
$colB[0,0] = 0.585894948492836
$colB[1,0] = 0.5443818955349968
$colB[2,0] = 0.5191546769189301
$colB[3,0] = 0.5089404865562699
$colB[4,0] = 0.5072484276242903
$colB[5,0] = 0.5190166568040127
$colB[6,0] = 0.5715270025814334
$colB[7,0] = 0.6765738234419985
$colB[8,0] = 0.7550190161279886
$colB[9,0] = 0.7909821150369396
$colB[10,0] = 0.8046402822871528
$colB[11,0] = 0.8016969889326049
$colB[12,0] = 0.7921049860224514
$colB[13,0] = 0.7862347728272994
$colB[14,0] = 0.752674315268365
$colB[15,0] = 0.7321570809357922
$colB[16,0] = 0.7203822908490736
$colB[17,0] = 0.7164000534678507
$colB[18,0] = 0.7343384678101188
$colB[19,0] = 0.7949213137420656
$colB[20,0] = 0.8347471213238009
$colB[21,0] = 0.8134702465914643
$colB[22,0] = 0.7333521978669921
$colB[23,0] = 0.6479335349488338
$ws.Range("B2:B25").Value = $colB

$colC = New-Object "object[,]" 24,1
$colC[0,0] = 0.05481939930005808
$colC[1,0] = 0.04840214020775591
$colC[2,0] = 0.04446839529647661
$colC[3,0] = 0.04286700172465885
$colC[4,0] = 0.0426011908830759
$colC[5,0] = 0.04444679168973664
$colC[6,0] = 0.05260538800790471
$colC[7,0] = 0.06865614875627557
$colC[8,0] = 0.08048195361878641
$colC[9,0] = 0.08586950535922711
$colC[10,0] = 0.08791078186575874
$colC[11,0] = 0.08747110650367063
$colC[12,0] = 0.08603741988258662
$colC[13,0] = 0.08515939188788479
$colC[14,0] = 0.0801300229466051
$colC[15,0] = 0.07704669046786705
$colC[16,0] = 0.07527398671498986
$colC[17,0] = 0.07467390848388789
$colC[18,0] = 0.07737483920418242
$colC[19,0] = 0.08645849789851923
$colC[20,0] = 0.09240176051008575
$colC[21,0] = 0.089229132852779
$colC[22,0] = 0.07722648333118798
$colC[23,0] = 0.06430827045370791
$ws.Range("C2:C25").Value = $colC

$colD = New-Object "object[,]" 24,1
$colD[0,0] = 0.2242817457523358
$colD[1,0] = 0.216537183565066
$colD[2,0] = 0.2118888394106904
$colD[3,0] = 0.2100215301519199
$colD[4,0] = 0.2097130931100963
$colD[5,0] = 0.2118635471103261
$colD[6,0] = 0.2215892778046111
$colD[7,0] = 0.2415077624960276
$colD[8,0] = 0.2566577326184074
$colD[9,0] = 0.2636619613549271
$colD[10,0] = 0.2663304165589011
$colD[11,0] = 0.2657550016536447
$colD[12,0] = 0.2638811743382803
$colD[13,0] = 0.2627354966271014
$colD[14,0] = 0.2562022469205942
$colD[15,0] = 0.2522230621713675
$colD[16,0] = 0.2499449265843907
$colD[17,0] = 0.2491754091006158
$colD[18,0] = 0.2526455579939864
$colD[19,0] = 0.264431126226043
$colD[20,0] = 0.2722275249224992
$colD[21,0] = 0.2680578733359198
$colD[22,0] = 0.2524545178802668
$colD[23,0] = 0.2360286884893412
$ws.Range("D2:D25").Value = $colD

$colF = New-Object "object[,]" 24,1
$colF[0,0] = 1.910637016552201
$colF[1,0] = 1.904055078749934
$colF[2,0] = 1.901035710885026
$colF[3,0] = 1.900062121247714
$colF[4,0] = 1.899915964721259
$colF[5,0] = 1.901021541040592
$colF[6,0] = 1.908155353643068
$colF[7,0] = 1.930263973362798
$colF[8,0] = 1.95147612376438
$colF[9,0] = 1.962209548121308
$colF[10,0] = 1.966430160636136
$colF[11,0] = 1.965514229992721
$colF[12,0] = 1.962553650681116
$colF[13,0] = 1.960760545842803
$colF[14,0] = 1.95079650011327
$colF[15,0] = 1.944961660944088
$colF[16,0] = 1.941707624482532
$colF[17,0] = 1.940623376709652
$colF[18,0] = 1.945572230334776
$colF[19,0] = 1.963419006276283
$colF[20,0] = 1.975992872680763
$colF[21,0] = 1.969198623924939
$colF[22,0] = 1.945295878900794
$colF[23,0] = 1.923411808069758
$ws.Range("F2:F25").Value = $colF

$colG = New-Object "object[,]" 24,1
$colG[0,0] = 1.174436044363546
$colG[1,0] = 1.171183168931677
$colG[2,0] = 1.169880721440464
$colG[3,0] = 1.169524435617362
$colG[4,0] = 1.169475805081802
$colG[5,0] = 1.169875210347385
$colG[6,0] = 1.173170086267746
$colG[7,0] = 1.185157747726777
$colG[8,0] = 1.197355391074325
$colG[9,0] = 1.203645380833564
$colG[10,0] = 1.206134160026764
$colG[11,0] = 1.205593397687025
$colG[12,0] = 1.203847990047464
$colG[13,0] = 1.202792807790601
$colG[14,0] = 1.196959267159613
$colG[15,0] = 1.193570643291068
$colG[16,0] = 1.191691345004216
$colG[17,0] = 1.191067016864196
$colG[18,0] = 1.193924147131199
$colG[19,0] = 1.20435775495811
$colG[20,0] = 1.211799941059411
$colG[21,0] = 1.207770778697679
$colG[22,0] = 1.193764113546735
$colG[23,0] = 1.181320936575545
$ws.Range("G2:G25").Value = $colG

$colH = New-Object "object[,]" 24,1
$colH[0,0] = 1.137536967554624
$colH[1,0] = 1.141032829470376
$colH[2,0] = 1.143677737759958
$colH[3,0] = 1.144880893753651
$colH[4,0] = 1.145088246786273
$colH[5,0] = 1.143693456473827
$colH[6,0] = 1.138638887855478
$colH[7,0] = 1.132682685225546
$colH[8,0] = 1.130720791658675
$colH[9,0] = 1.130353156005853
$colH[10,0] = 1.130289452861788
$colH[11,0] = 1.130299813359514
$colH[12,0] = 1.130346401485113
$colH[13,0] = 1.130384773265078
$colH[14,0] = 1.130755382111005
$colH[15,0] = 1.131117191429311
$colH[16,0] = 1.131374692376951
$colH[17,0] = 1.131470360671031
$colH[18,0] = 1.131073563424778
$colH[19,0] = 1.130330667698786
$colH[20,0] = 1.130285307372787
$colH[21,0] = 1.130269227856132
$colH[22,0] = 1.131093133476099
$colH[23,0] = 1.133870241899615
$ws.Range("H2:H25").Value = $colH

$colJ = New-Object "object[,]" 24,1
$colJ[0,0] = 0.3517343430940798
$colJ[1,0] = 0.3406371280884315
$colJ[2,0] = 0.3340338975760204
$colJ[3,0] = 0.331395858851252
$colJ[4,0] = 0.3309610024711276
$colJ[5,0] = 0.3339981062834454
$colJ[6,0] = 0.3478642675218424
$colJ[7,0] = 0.3767324746655305
$colJ[8,0] = 0.3989757256897377
$colJ[9,0] = 0.4093219029795421
$colJ[10,0] = 0.4132726207077155
$colJ[11,0] = 0.4124203003361373
$colJ[12,0] = 0.4096462719835614
$colJ[13,0] = 0.40795138031649
$colJ[14,0] = 0.3983041685789885
$colJ[15,0] = 0.392444294691046
$colJ[16,0] = 0.3890952671336834
$colJ[17,0] = 0.3879650188349757
$colJ[18,0] = 0.3930658708729027
$colJ[19,0] = 0.4104601790339188
$colJ[20,0] = 0.4220198641047688
$colJ[21,0] = 0.4158326823346812
$colJ[22,0] = 0.3927847943328544
$colJ[23,0] = 0.3687421587036397
$ws.Range("J2:J25").Value = $colJ

$colK = New-Object "object[,]" 24,1
$colK[0,0] = 0.5636292628394983
$colK[1,0] = 0.5183171290006499
$colK[2,0] = 0.4907217174864797
$colK[3,0] = 0.4795335094514996
$colK[4,0] = 0.4776791738239012
$colK[5,0] = 0.4905705977363652
$colK[6,0] = 0.5479587559967456
$colK[7,0] = 0.6622903425867435
$colK[8,0] = 0.7473884728138955
$colK[9,0] = 0.7863421293524766
$colK[10,0] = 0.8011276144369504
$colK[11,0] = 0.7979417613547071
$colK[12,0] = 0.7875578483782419
$colK[13,0] = 0.7812018984486144
$colK[14,0] = 0.7448476113870299
$colK[15,0] = 0.7226072735617493
$colK[16,0] = 0.7098380198815164
$colK[17,0] = 0.7055184930694622
$colK[18,0] = 0.7249724374190407
$colK[19,0] = 0.7906069186529407
$colK[20,0] = 0.8337043034274245
$colK[21,0] = 0.8106840607281924
$colK[22,0] = 0.7239030938305575
$colK[23,0] = 0.6311680639998656
$ws.Range("K2:K25").Value = $colK

$colN = New-Object "object[,]" 24,1
$colN[0,0] = 1.918041824904009
$colN[1,0] = 1.937447975651221
$colN[2,0] = 1.949999629008044
$colN[3,0] = 1.955274397846221
$colN[4,0] = 1.956159921695708
$colN[5,0] = 1.950070119235374
$colN[6,0] = 1.924600958727449
$colN[7,0] = 1.879708130977448
$colN[8,0] = 1.849809975982371
$colN[9,0] = 1.836879357042129
$colN[10,0] = 1.832079373321875
$colN[11,0] = 1.833108838102845
$colN[12,0] = 1.836482523182319
$colN[13,0] = 1.838561583414403
$colN[14,0] = 1.850668524396564
$colN[15,0] = 1.858267527028296
$colN[16,0] = 1.862701343474569
$colN[17,0] = 1.864213384444312
$colN[18,0] = 1.857452071252197
$colN[19,0] = 1.835488967947683
$colN[20,0] = 1.821697731369717
$colN[21,0] = 1.829006803848152
$colN[22,0] = 1.85782053618632
$colN[23,0] = 1.891311495463974
$ws.Range("N2:N25").Value = $colN

